$d = $word.ActiveDocument

# --- Remove the existing "_GoBack" bookmark from the end of the intro paragraph.
# It will be re-created later, right after "Database Design" in the title.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Title: "Database Documentation" -> "Database Design Documentation"
$d.Content.Find.Execute("Database Documentation", $false, $false, $false, $false, $false, $true, 1, $false, "Database Design Documentation", 2) | Out-Null

# --- Body paragraph fixes
# "the devices name and" -> "the devices' names and" (curly apostrophe)
$apos = [char]8217
$d.Content.Find.Execute("the devices name and", $false, $false, $false, $false, $false, $true, 1, $false, ("the devices" + $apos + " names and"), 2) | Out-Null

# "...was first online and the date and date the device was last online..." ->
# "...was first online and the time and date the device was last online..."
$d.Content.Find.Execute("was first online and the date and date the device was last", $false, $false, $false, $false, $false, $true, 1, $false, "was first online and the time and date the device was last", 2) | Out-Null

# --- Re-insert the "_GoBack" bookmark right after "Database Design" in the title.
$titleRange = $d.Paragraphs(1).Range
$pos = $titleRange.Start + ("Database Design").Length
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
